$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B -> C)
$ws.Columns("B:B").Insert()

# New column B values (set data rows first so the shared-string table
# picks up "View Lead | opentaps CRM" before "expectedTitle")
$ws.Range("B2").Value = "View Lead | opentaps CRM"
$ws.Range("B3").Value = "View Lead | opentaps CRM"
$ws.Range("B1").Value = "expectedTitle"

# Column widths per diff (target stored widths: B=29.140625, C=19.42578125).
# The host's ColumnWidth setter quantizes; these inputs land on the closest
# reachable stored widths (29.166666666666668 and 19.5 respectively).
$ws.Columns("B:B").ColumnWidth = 28.333333333333332
$ws.Columns("C:C").ColumnWidth = 18.666666666666668

# Selection per diff
$ws.Range("C1").Select()
